# Minor update from Alison on geometry mapping reference file.
#
# The "Name PDB"-adjacent "Name FindGeo" column (column C) contained several
# labels that used a comma after the base-polyhedron name, e.g.
#   "octahedron, face monocapped ..."
#   "trigonal prism, square-face ..."
#   "square antiprism, square-face ..."
# Alison's edit drops that stray comma from each of those labels. Apply the
# same text fix-up to the corresponding cells on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C26").Value = "octahedron face monocapped with a vacancy (capped face)"
$ws.Range("C27").Value = "octahedron face monocapped with a vacancy (non-capped face)"
$ws.Range("C28").Value = "trigonal prism square-face monocapped with a vacancy (capped face)"
$ws.Range("C63").Value = "square antiprism square-face monocapped"
$ws.Range("C62").Value = "trigonal prism square-face tricapped"
$ws.Range("C51").Value = "trigonal prism square-face bicapped"
$ws.Range("C52").Value = "trigonal prism triangular-face bicapped"
$ws.Range("C50").Value = "octahedron trans-bicapped"
$ws.Range("C36").Value = "trigonal prism square-face monocapped"
$ws.Range("C37").Value = "octahedron face monocapped"
$ws.Range("C29").Value = "trigonal prism square-face monocapped with a vacancy (non-capped edge)`t"

# Reflect the author's final cursor position/selection in the saved view.
$ws.Range("C30").Select()
